$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column D ("Tipo") to make room for "MAE"
$ws.Range("D1").EntireColumn.Insert()

# New header for the inserted column D
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "MAE"

# Updated values
$ws.Range("B2").Value = 0.04908998251238541
$ws.Range("C2").Value = 0.9985563473642154
$ws.Range("D2").Value = 0.1617162528759863
